# Fix a couple of shared-string typos/content and the view state
# (zoom level + active selection), matching the upstream template update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Typo fix: "Infirimière" -> "Infirmière"
$ws.Range("B6").Value = "Infirmière"

# Icon swap for "Activités physiques": weight-lifter -> runner
$ws.Range("C15").Value = "🏃"

# Move the selection to C15 and bump the zoom level to 160%
$null = $ws.Range("C15").Select()
$excel.ActiveWindow.Zoom = 160
